$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5. This shifts the existing rows 5-53 down to 6-54,
# preserving all of their data/formatting untouched.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44882
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 100112042
$ws.Range("G5").Value = "Locoto"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("N5").Value = "$/kilo"
$ws.Range("O5").Value = "Región de Arica y Parinacota"
$ws.Range("P5").Value = 2500
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = "Hortaliza"
